# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets share identical data, and the same set of row updates applies to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 45
    5  = 246
    6  = 30
    7  = 108
    8  = 253
    13 = 79
    14 = 322
    16 = 462
    17 = 378
    18 = 133
    19 = 58
    20 = 29
    22 = 896
    23 = 2680
    24 = 22
    26 = 519
    27 = 960
    28 = 565
    29 = 448
    30 = 256
    33 = 590
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
